$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data update (automatic electricity price refresh)
$ws.Range("A2").Value = 46019
$ws.Range("B2").Value = 62.5
$ws.Range("C2").Value = 63.75
$ws.Range("D2").Value = 54.51
$ws.Range("E2").Value = 48.86
$ws.Range("F2").Value = 51.86
$ws.Range("G2").Value = 55.2
$ws.Range("H2").Value = 65.87
$ws.Range("I2").Value = 70.31
$ws.Range("J2").Value = 79.54000000000001
$ws.Range("K2").Value = 83
$ws.Range("L2").Value = 79.09999999999999
$ws.Range("M2").Value = 72.37
$ws.Range("N2").Value = 79.31
$ws.Range("O2").Value = 76.15000000000001
$ws.Range("P2").Value = 78.09
$ws.Range("Q2").Value = 82.43000000000001
$ws.Range("R2").Value = 91.44
$ws.Range("S2").Value = 98.56
$ws.Range("T2").Value = 103.05
$ws.Range("U2").Value = 105.21
$ws.Range("V2").Value = 105.01
$ws.Range("W2").Value = 99.47
$ws.Range("X2").Value = 98.28
$ws.Range("Y2").Value = 90.63
$ws.Range("Z2").Value = 78.94
$ws.Range("AB2").Value = 99.56
$ws.Range("AD2").Value = 104.13
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 102.24
$ws.Range("AG2").Value = "0h-14h"
